$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (report header) ---
$ws.Range("A8").Value = "Volume 31   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/11/2024  Through  11/17/2024"

# --- Crime statistics table updates (rows 16-31) ---
# Row 16
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 114
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = -12.307692307692
$ws.Range("L16").Value = 17.525773195876
$ws.Range("M16").Value = -22.448979591836
$ws.Range("N16").Value = -83.620689655172

# Row 17
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -7.692307692307
$ws.Range("I17").Value = 160
$ws.Range("J17").Value = 163
$ws.Range("K17").Value = -1.840490797546
$ws.Range("L17").Value = 23.076923076923
$ws.Range("M17").Value = 61.616161616161
$ws.Range("N17").Value = -40.298507462686

# Row 18
$ws.Range("C18").Value = "'0"
$ws.Range("C20").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = -44.047619047619
$ws.Range("L18").Value = -57.272727272727
$ws.Range("M18").Value = -33.802816901408
$ws.Range("N18").Value = -91.081593927893

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -27.586206896551
$ws.Range("I19").Value = 306
$ws.Range("J19").Value = 348
$ws.Range("K19").Value = -12.068965517241
$ws.Range("L19").Value = -0.970873786407
$ws.Range("M19").Value = 25.925925925925
$ws.Range("N19").Value = -49.169435215946

# Row 20
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -16.666666666666
$ws.Range("L20").Value = -42.647058823529
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -88.145896656535

# Row 21
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -47.058823529411
$ws.Range("F21").Value = 52
$ws.Range("G21").Value = 69
$ws.Range("H21").Value = -24.63768115942
$ws.Range("I21").Value = 673
$ws.Range("J21").Value = 831
$ws.Range("K21").Value = -19.013237063778
$ws.Range("L21").Value = -7.808219178082
$ws.Range("M21").Value = 12.166666666666
$ws.Range("N21").Value = -72.608872608872

# Row 22
$ws.Range("C22").Value = "'0"
$ws.Range("C20").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("E16").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("J22").Value = 33
$ws.Range("K22").Value = -30.30303030303
$ws.Range("M22").Value = -23.333333333333

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = -33.333333333333
$ws.Range("I23").Value = 118
$ws.Range("J23").Value = 137
$ws.Range("K23").Value = -13.868613138686
$ws.Range("L23").Value = 15.686274509803
$ws.Range("M23").Value = 68.571428571428

# Row 24
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 37
$ws.Range("H24").Value = 2.702702702702
$ws.Range("I24").Value = 372
$ws.Range("J24").Value = 451
$ws.Range("K24").Value = -17.516629711751
$ws.Range("L24").Value = -18.06167400881
$ws.Range("M24").Value = -32.116788321167

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 13
$ws.Range("H25").Value = 85.714285714285
$ws.Range("I25").Value = 70
$ws.Range("J25").Value = 145
$ws.Range("K25").Value = -51.724137931034
$ws.Range("L25").Value = -51.048951048951

# Row 26
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = -13.043478260869
$ws.Range("I26").Value = 216
$ws.Range("J26").Value = 214
$ws.Range("K26").Value = 0.934579439252
$ws.Range("L26").Value = 17.391304347826
$ws.Range("M26").Value = -18.490566037735

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 37
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = -13.953488372093
$ws.Range("L28").Value = -24.489795918367

# Row 29
$ws.Range("N29").Value = -77.419354838709

# Row 30
$ws.Range("N30").Value = -76.923076923076

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = -80
$ws.Range("J31").Value = 12
$ws.Range("K31").Value = -25
